$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Flip the "more than 1 rivers" marker cells that are no longer single matches back to 0,
# and mark the newly-found multi-river cells (row 6) as 1.
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 1

$ws.Range("F7").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 0

$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0

$ws.Range("G17").Value = 0

$ws.Range("G18").Value = 0

$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0

$ws.Range("F20").Value = 0

$ws.Range("F21").Value = 0

$ws.Range("F22").Value = 0

$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0

$ws.Range("G24").Value = 0

$ws.Range("E25").Value = 0
$ws.Range("H25").Value = 0

$ws.Range("F26").Value = 0
$ws.Range("I26").Value = 0

$ws.Range("G27").Value = 0
$ws.Range("I27").Value = 0

$ws.Range("G28").Value = 0
$ws.Range("I28").Value = 0

$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0

# New search-engine-style summary cell: count how many cells in the grid equal 1.
$ws.Range("G34").Formula = "=COUNTIF(A1:O30,""=1"")"

# View state: zoom in and move the active selection.
$ws.Application.ActiveWindow.Zoom = 70
$ws.Range("K32").Select()
